$wb = $excel.ActiveWorkbook

# Select the whole sheet on Slovakia (this is what's reflected in the saved
# view state once it stops being the active tab)
$slovakia = $wb.Worksheets.Item("Slovakia")
[void]$slovakia.Range("A1:XFD1048576").Select()

# Duplicate the Slovakia sheet and place the copy right after it, at the end
$slovakia.Copy($null, $slovakia)
$hungary = $wb.ActiveSheet
$hungary.Name = "Hungary"

# Update the market-specific content for Hungary
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3593/T3618"

# Match the saved selection/active cell state and make Hungary the active tab
[void]$hungary.Range("B2:B4").Select()
$hungary.Activate()
